$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 45 (shifting existing
# rows 45-76 down to 47-78). This matches the dimension change
# A1:T76 -> A1:T78 seen in the diff.
$ws.Rows.Item(45).Insert()
$ws.Rows.Item(45).Insert()

# The rows that used to be 45 and 46 are now at 47 and 48. Duplicate their
# full formatting/content into the two freshly inserted blank rows (45/46)
# so every column besides Fecha (D) and Volumen (M) matches.
$ws.Rows.Item(47).Copy()
$ws.Rows.Item(45).PasteSpecial()
$ws.Rows.Item(48).Copy()
$ws.Rows.Item(46).PasteSpecial()

# New row 45: Fecha 2022-02-10 (serial 44587), Volumen 180
$ws.Cells.Item(45, 4).Value = 44587
$ws.Cells.Item(45, 13).Value = 180

# New row 46: Fecha 2022-02-10 (serial 44587), Volumen 260
$ws.Cells.Item(46, 4).Value = 44587
$ws.Cells.Item(46, 13).Value = 260
